# Update cryptocurrency price/volume figures (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.604.81'
$ws.Range("E2").Value = '  -5.98%  '
$ws.Range("D3").Value = '3.263.19'
$ws.Range("E3").Value = '  -6.77%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.05'
$ws.Range("E5").Value = '  -4.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.10'
$ws.Range("E6").Value = '  -4.80%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("E8").Value = '  -4.34%  '
$ws.Range("D9").Value = '3.254.59'
$ws.Range("E9").Value = '  -6.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.182'
$ws.Range("E10").Value = '  -10.75%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.580'
$ws.Range("E11").Value = '  -6.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.98'
$ws.Range("E12").Value = '  -9.08%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000265'
$ws.Range("E13").Value = '  -7.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.58'
$ws.Range("E14").Value = '  -6.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '629.11'
$ws.Range("E15").Value = '  -2.79%  '
$ws.Range("D16").Value = '3.792.76'
$ws.Range("E16").Value = '  -6.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '17.98'
$ws.Range("E17").Value = '  -1.50%  '
$ws.Range("D18").Value = '65.634.98'
$ws.Range("E19").Value = '  -3.50%  '
$ws.Range("D20").Value = '3.270.55'
$ws.Range("E20").Value = '  -6.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.25'
$ws.Range("E21").Value = '  -8.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.898'
$ws.Range("E22").Value = '  -5.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '18.16'
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '106.70'
$ws.Range("E24").Value = '  +7.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.88'
$ws.Range("E25").Value = '  -6.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.95'
$ws.Range("E26").Value = '  -7.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.66'
$ws.Range("E27").Value = '  -7.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.49'
$ws.Range("E28").Value = '  -5.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.63'
$ws.Range("E29").Value = '  -7.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.16'
$ws.Range("E30").Value = '  -7.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.91'
$ws.Range("E31").Value = '  -8.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.21'
$ws.Range("E32").Value = '  -7.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.99'
$ws.Range("E33").Value = '  -5.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.104'
$ws.Range("E34").Value = '  -5.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.59'
$ws.Range("E35").Value = '  -6.07%  '
$ws.Range("D36").Value = '3.716.19'
$ws.Range("E36").Value = '  +0.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '520.68'
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.37'
$ws.Range("E39").Value = '  -4.78%  '
$ws.Range("D40").Value = '0.0₃0729'
$ws.Range("E40").Value = '  -7.62%  '
$ws.Range("E41").Value = '  -1.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.69'
$ws.Range("E42").Value = '  -7.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.36'
$ws.Range("E43").Value = '  -5.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '32.77'
$ws.Range("E44").Value = '  -4.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.335'
$ws.Range("E45").Value = '  -10.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.24'
$ws.Range("E46").Value = '  -4.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0412'
$ws.Range("E47").Value = '  -7.09%  '
$ws.Range("E48").Value = '  -4.44%  '
$ws.Range("E49").Value = '  -8.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  +0.23%  '
$ws.Range("E51").Value = '  +0.45%  '
